$wb = $excel.ActiveWorkbook

# This workbook scrapes bilibili-show event listings for Suzhou into four
# sheets; "展览" and "全部类型" carry the same 14-row exhibition table. The
# refresh re-scrapes the "想去人数" (want-to-go) counters, inserts a newly
# discovered event ("苏州·绘时国乙1.0-秩序之外", 2024.04.13) as the new
# row 12 (sorted by date before the 04.21 event), and pushes the rest of
# the table down by one row.

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Refresh the "想去人数" (want-to-go) counts for the rows that keep
    # their position in the table.
    $ws.Range("F2").Value = 1827
    $ws.Range("F3").Value = 252
    $ws.Range("F4").Value = 239
    $ws.Range("F5").Value = 7941
    $ws.Range("F6").Value = 569
    $ws.Range("F7").Value = 559
    $ws.Range("F8").Value = 83
    $ws.Range("F9").Value = 21
    $ws.Range("F10").Value = 9148
    $ws.Range("F11").Value = 2385

    # Insert the newly scraped event as row 12; this shifts the former
    # rows 12-14 down to 13-15 and grows the sheet dimension to A1:J15.
    $ws.Rows.Item(12).Insert()

    # The inserted row's A cell needs the same bold/centered/thin-border
    # style used by every other row-number cell in column A.
    $ws.Range("A12").Borders.LineStyle = 1
    $ws.Range("A12").Font.Bold = $true
    $ws.Range("A12").HorizontalAlignment = -4108
    $ws.Range("A12").VerticalAlignment = -4160

    $ws.Range("A12").Value = 11
    $ws.Range("B12").Value = "'2024.04.13"
    $ws.Range("C12").Value = "苏州·绘时国乙1.0-秩序之外"
    $ws.Range("D12").Value = "石路步行街永福桥浜15号 银河广场"
    $ws.Range("E12").Value = "2024.04.13 13:30-04.13 20:00"
    $ws.Range("F12").Value = 7
    $ws.Range("G12").Value = "'78"
    $ws.Range("H12").Value = $false
    $ws.Range("I12").Value = "https://show.bilibili.com/platform/detail.html?id=80789&msource=Msearch_colligation"
    $ws.Range("J12").Value = "//i0.hdslb.com/bfs/openplatform/202401/SjKfDxBh1705041298410.jpeg"

    # The rows that were pushed down (formerly 12/13/14, now 13/14/15)
    # keep all of their original text, but their want-to-go counters were
    # refreshed too.
    $ws.Range("F13").Value = 306
    $ws.Range("F14").Value = 9959
    $ws.Range("F15").Value = 10531
}
